$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "261.18"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.00%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.98"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.41%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.701"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.28%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06218"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.46%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.751"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.37%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8519"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.18%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9161"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.66%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1402"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.37%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04857"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.38%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07090"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.22%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03109"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.86%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09066"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.18%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001546"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.90%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006176"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.86%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006016"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.58%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.442"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.43%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.176"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.53%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.47%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1311"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.98%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.092"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.50%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04234"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.35%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001203"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.33%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.00%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "4.39%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03951"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.79%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.11%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004118"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.45%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.07%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-7.54%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005163"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.41%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.04%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2103"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "59.18%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.04%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.04%"
